{"js": "// Replace each two-digit multiplication prompt in the worksheet with its\n// new value, in document order. Every \"old\" value below is unique in the\n// document and none of the \"new\" values collide with any \"old\" value, so\n// doing the replacements as independent search+replace operations (rather\n// than relying on positional indices) is safe and order-independent.\nconst replacements = [\n  [\"66\u00d711=\", \"69\u00d726=\"],\n  [\"32\u00d740=\", \"98\u00d768=\"],\n  [\"80\u00d794=\", \"23\u00d751=\"],\n  [\"90\u00d756=\", \"36\u00d715=\"],\n  [\"83\u00d792=\", \"26\u00d775=\"],\n  [\"60\u00d772=\", \"97\u00d733=\"],\n  [\"31\u00d799=\", \"75\u00d731=\"],\n  [\"44\u00d723=\", \"92\u00d713=\"],\n  [\"89\u00d765=\", \"73\u00d768=\"],\n  [\"46\u00d765=\", \"16\u00d753=\"],\n  [\"86\u00d778=\", \"72\u00d742=\"],\n  [\"46\u00d741=\", \"64\u00d750=\"],\n  [\"26\u00d791=\", \"70\u00d781=\"],\n  [\"95\u00d754=\", \"45\u00d738=\"],\n  [\"58\u00d731=\", \"27\u00d797=\"],\n  [\"70\u00d714=\", \"46\u00d737=\"],\n  [\"70\u00d743=\", \"47\u00d756=\"],\n  [\"59\u00d735=\", \"19\u00d715=\"],\n  [\"96\u00d765=\", \"18\u00d750=\"],\n  [\"62\u00d740=\", \"38\u00d745=\"],\n  [\"22\u00d796=\", \"50\u00d771=\"],\n  [\"58\u00d750=\", \"16\u00d785=\"],\n  [\"22\u00d798=\", \"49\u00d740=\"],\n  [\"49\u00d789=\", \"41\u00d725=\"],\n  [\"18\u00d748=\", \"27\u00d737=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication prompt in the worksheet with its\n# new value. Every \"old\" value is unique in the document and none of the\n# \"new\" values collide with any \"old\" value, so running these as\n# independent Find/Replace passes over the whole document is safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"66\u00d711=\", \"69\u00d726=\"),\n    @(\"32\u00d740=\", \"98\u00d768=\"),\n    @(\"80\u00d794=\", \"23\u00d751=\"),\n    @(\"90\u00d756=\", \"36\u00d715=\"),\n    @(\"83\u00d792=\", \"26\u00d775=\"),\n    @(\"60\u00d772=\", \"97\u00d733=\"),\n    @(\"31\u00d799=\", \"75\u00d731=\"),\n    @(\"44\u00d723=\", \"92\u00d713=\"),\n    @(\"89\u00d765=\", \"73\u00d768=\"),\n    @(\"46\u00d765=\", \"16\u00d753=\"),\n    @(\"86\u00d778=\", \"72\u00d742=\"),\n    @(\"46\u00d741=\", \"64\u00d750=\"),\n    @(\"26\u00d791=\", \"70\u00d781=\"),\n    @(\"95\u00d754=\", \"45\u00d738=\"),\n    @(\"58\u00d731=\", \"27\u00d797=\"),\n    @(\"70\u00d714=\", \"46\u00d737=\"),\n    @(\"70\u00d743=\", \"47\u00d756=\"),\n    @(\"59\u00d735=\", \"19\u00d715=\"),\n    @(\"96\u00d765=\", \"18\u00d750=\"),\n    @(\"62\u00d740=\", \"38\u00d745=\"),\n    @(\"22\u00d796=\", \"50\u00d771=\"),\n    @(\"58\u00d750=\", \"16\u00d785=\"),\n    @(\"22\u00d798=\", \"49\u00d740=\"),\n    @(\"49\u00d789=\", \"41\u00d725=\"),\n    @(\"18\u00d748=\", \"27\u00d737=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
